# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# Commit: Updated cryptos list on Sat Apr 15 06:48:50 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.468.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.091.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5211"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.65%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08925"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.79%  "

$ws.Range("E12").Value = "  -4.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.093.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.678"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.689"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("E16").Value = "  -2.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.254"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.517.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.307"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.341.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.26%  "

$ws.Range("E27").Value = "  -3.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.558"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.188"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1068"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.650"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.149"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.898"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02555"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06803"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.465"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2249"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6862"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.248"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6315"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.28%  "

$ws.Range("E47").Value = "  -3.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.627"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.237"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.245"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.23%  "
